$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 12.47064466666667
$ws.Range("H2").Value2 = 37.411934
$ws.Range("I2").Value2 = 0.183409848855644
$ws.Range("J2").Value2 = 0.183409848855644
$ws.Range("M2").Value2 = 17.10933733333333
$ws.Range("N2").Value2 = 51.328012
$ws.Range("O2").Value2 = 0.3554368716515803
$ws.Range("P2").Value2 = 0.3554368716515803
$ws.Range("Q2").Value2 = 213.3644663661342
$ws.Range("R2").Value2 = 1920.280197295208
$ws.Range("S2").Value2 = 0.06519062290733928
$ws.Range("T2").Value2 = 0.06519062290733928
$ws.Range("G3").Value2 = 12.47064466666667
$ws.Range("H3").Value2 = 37.411934
$ws.Range("I3").Value2 = 0.183409848855644
$ws.Range("J3").Value2 = 0.183409848855644
$ws.Range("O3").Value2 = 0.2270123898818874
$ws.Range("P3").Value2 = 0.2270123898818874
$ws.Range("Q3").Value2 = 136.2727991628562
$ws.Range("R3").Value2 = 1226.455192465706
$ws.Range("S3").Value2 = 0.0416363081165955
$ws.Range("T3").Value2 = 0.04163630811659549
$ws.Range("G4").Value2 = 12.47064466666667
$ws.Range("H4").Value2 = 37.411934
$ws.Range("I4").Value2 = 0.183409848855644
$ws.Range("J4").Value2 = 0.183409848855644
$ws.Range("M4").Value2 = 11.616679
$ws.Range("N4").Value2 = 34.850037
$ws.Range("O4").Value2 = 0.2413299803667016
$ws.Range("P4").Value2 = 0.2413299803667016
$ws.Range("Q4").Value2 = 144.8674760157287
$ws.Range("R4").Value2 = 1303.807284141558
$ws.Range("S4").Value2 = 0.04426229522339228
$ws.Range("T4").Value2 = 0.04426229522339228
$ws.Range("G5").Value2 = 12.47064466666667
$ws.Range("H5").Value2 = 37.411934
$ws.Range("I5").Value2 = 0.183409848855644
$ws.Range("J5").Value2 = 0.183409848855644
$ws.Range("M5").Value2 = 8.482576333333332
$ws.Range("N5").Value2 = 25.447729
$ws.Range("O5").Value2 = 0.1762207580998305
$ws.Range("P5").Value2 = 0.1762207580998305
$ws.Range("Q5").Value2 = 105.7831953108762
$ws.Range("R5").Value2 = 952.048757797886
$ws.Range("S5").Value2 = 0.03232062260831692
$ws.Range("T5").Value2 = 0.03232062260831692
$ws.Range("I6").Value2 = 0.1770741628042856
$ws.Range("J6").Value2 = 0.1770741628042856
$ws.Range("M6").Value2 = 17.10933733333333
$ws.Range("N6").Value2 = 51.328012
$ws.Range("O6").Value2 = 0.3554368716515803
$ws.Range("P6").Value2 = 0.3554368716515803
$ws.Range("Q6").Value2 = 205.9940318892191
$ws.Range("R6").Value2 = 1853.946287002972
$ws.Range("S6").Value2 = 0.06293868647747791
$ws.Range("T6").Value2 = 0.06293868647747791
$ws.Range("I7").Value2 = 0.1770741628042856
$ws.Range("J7").Value2 = 0.1770741628042856
$ws.Range("O7").Value2 = 0.2270123898818874
$ws.Range("P7").Value2 = 0.2270123898818874
$ws.Range("S7").Value2 = 0.0401980288845353
$ws.Range("T7").Value2 = 0.0401980288845353
$ws.Range("I8").Value2 = 0.1770741628042856
$ws.Range("J8").Value2 = 0.1770741628042856
$ws.Range("M8").Value2 = 11.616679
$ws.Range("N8").Value2 = 34.850037
$ws.Range("O8").Value2 = 0.2413299803667016
$ws.Range("P8").Value2 = 0.2413299803667016
$ws.Range("Q8").Value2 = 139.8631926971663
$ws.Range("R8").Value2 = 1258.768734274497
$ws.Range("S8").Value2 = 0.04273330423300838
$ws.Range("T8").Value2 = 0.04273330423300838
$ws.Range("I9").Value2 = 0.1770741628042856
$ws.Range("J9").Value2 = 0.1770741628042856
$ws.Range("M9").Value2 = 8.482576333333332
$ws.Range("N9").Value2 = 25.447729
$ws.Range("O9").Value2 = 0.1762207580998305
$ws.Range("P9").Value2 = 0.1762207580998305
$ws.Range("Q9").Value2 = 102.1290343201721
$ws.Range("R9").Value2 = 919.1613088815488
$ws.Range("S9").Value2 = 0.03120414320926403
$ws.Range("T9").Value2 = 0.03120414320926403
$ws.Range("G10").Value2 = 7.218786333333333
$ws.Range("H10").Value2 = 21.656359
$ws.Range("I10").Value2 = 0.1061690510561032
$ws.Range("J10").Value2 = 0.1061690510561032
$ws.Range("M10").Value2 = 17.10933733333333
$ws.Range("N10").Value2 = 51.328012
$ws.Range("O10").Value2 = 0.3554368716515803
$ws.Range("P10").Value2 = 0.3554368716515803
$ws.Range("Q10").Value2 = 123.5086505142564
$ws.Range("R10").Value2 = 1111.577854628308
$ws.Range("S10").Value2 = 0.03773639537359825
$ws.Range("T10").Value2 = 0.03773639537359824
$ws.Range("G11").Value2 = 7.218786333333333
$ws.Range("H11").Value2 = 21.656359
$ws.Range("I11").Value2 = 0.1061690510561032
$ws.Range("J11").Value2 = 0.1061690510561032
$ws.Range("O11").Value2 = 0.2270123898818874
$ws.Range("P11").Value2 = 0.2270123898818874
$ws.Range("Q11").Value2 = 78.88318900075345
$ws.Range("R11").Value2 = 709.948701006781
$ws.Range("S11").Value2 = 0.02410169001173812
$ws.Range("T11").Value2 = 0.02410169001173812
$ws.Range("G12").Value2 = 7.218786333333333
$ws.Range("H12").Value2 = 21.656359
$ws.Range("I12").Value2 = 0.1061690510561032
$ws.Range("J12").Value2 = 0.1061690510561032
$ws.Range("M12").Value2 = 11.616679
$ws.Range("N12").Value2 = 34.850037
$ws.Range("O12").Value2 = 0.2413299803667016
$ws.Range("P12").Value2 = 0.2413299803667016
$ws.Range("Q12").Value2 = 83.85832360392033
$ws.Range("R12").Value2 = 754.7249124352829
$ws.Range("S12").Value2 = 0.02562177500692074
$ws.Range("T12").Value2 = 0.02562177500692074
$ws.Range("G13").Value2 = 7.218786333333333
$ws.Range("H13").Value2 = 21.656359
$ws.Range("I13").Value2 = 0.1061690510561032
$ws.Range("J13").Value2 = 0.1061690510561032
$ws.Range("M13").Value2 = 8.482576333333332
$ws.Range("N13").Value2 = 25.447729
$ws.Range("O13").Value2 = 0.1762207580998305
$ws.Range("P13").Value2 = 0.1762207580998305
$ws.Range("Q13").Value2 = 61.23390610652343
$ws.Range("R13").Value2 = 551.105154958711
$ws.Range("S13").Value2 = 0.01870919066384612
$ws.Range("T13").Value2 = 0.01870919066384612
$ws.Range("G14").Value2 = 36.26402933333333
$ws.Range("H14").Value2 = 108.792088
$ws.Range("I14").Value2 = 0.5333469372839672
$ws.Range("J14").Value2 = 0.5333469372839672
$ws.Range("M14").Value2 = 17.10933733333333
$ws.Range("N14").Value2 = 51.328012
$ws.Range("O14").Value2 = 0.3554368716515803
$ws.Range("P14").Value2 = 0.3554368716515803
$ws.Range("Q14").Value2 = 620.4535109298951
$ws.Range("R14").Value2 = 5584.081598369056
$ws.Range("S14").Value2 = 0.1895711668931649
$ws.Range("T14").Value2 = 0.1895711668931649
$ws.Range("G15").Value2 = 36.26402933333333
$ws.Range("H15").Value2 = 108.792088
$ws.Range("I15").Value2 = 0.5333469372839672
$ws.Range("J15").Value2 = 0.5333469372839672
$ws.Range("O15").Value2 = 0.2270123898818874
$ws.Range("P15").Value2 = 0.2270123898818874
$ws.Range("Q15").Value2 = 396.2746849315992
$ws.Range("R15").Value2 = 3566.472164384393
$ws.Range("S15").Value2 = 0.1210763628690185
$ws.Range("T15").Value2 = 0.1210763628690185
$ws.Range("G16").Value2 = 36.26402933333333
$ws.Range("H16").Value2 = 108.792088
$ws.Range("I16").Value2 = 0.5333469372839672
$ws.Range("J16").Value2 = 0.5333469372839672
$ws.Range("M16").Value2 = 11.616679
$ws.Range("N16").Value2 = 34.850037
$ws.Range("O16").Value2 = 0.2413299803667016
$ws.Range("P16").Value2 = 0.2413299803667016
$ws.Range("Q16").Value2 = 421.2675880119173
$ws.Range("R16").Value2 = 3791.408292107256
$ws.Range("S16").Value2 = 0.1287126059033803
$ws.Range("T16").Value2 = 0.1287126059033803
$ws.Range("G17").Value2 = 36.26402933333333
$ws.Range("H17").Value2 = 108.792088
$ws.Range("I17").Value2 = 0.5333469372839672
$ws.Range("J17").Value2 = 0.5333469372839672
$ws.Range("M17").Value2 = 8.482576333333332
$ws.Range("N17").Value2 = 25.447729
$ws.Range("O17").Value2 = 0.1762207580998305
$ws.Range("P17").Value2 = 0.1762207580998305
$ws.Range("Q17").Value2 = 307.6123969742391
$ws.Range("R17").Value2 = 2768.511572768152
$ws.Range("S17").Value2 = 0.09398680161840346
$ws.Range("T17").Value2 = 0.09398680161840348
